$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.337401628494263
$ws.Range("B1").Value = 3.584694385528564
$ws.Range("C1").Value = 5.460087776184082
$ws.Range("D1").Value = 1.415757775306702
$ws.Range("E1").Value = 0.6786640882492065
